$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1389.091
$ws.Range("I2").Value = 1171.5714
$ws.Range("J2").Value = 1769.75
$ws.Range("K2").Value = 1171.5714
$ws.Range("L2").Value = 1769.75
$ws.Range("M2").Value = -1058.5714
$ws.Range("N2").Value = -1995.75

$ws.Range("H15").Value = 1003.4048
$ws.Range("I15").Value = 1003.4048
$ws.Range("K15").Value = 3010.2144
$ws.Range("M15").Value = -2841.2144

$ws.Range("H28").Value = 653.0833
$ws.Range("I28").Value = 612
$ws.Range("K28").Value = 612
$ws.Range("M28").Value = -127

$ws.Range("H32").Value = 11044
$ws.Range("J32").Value = 3933.3333
$ws.Range("L32").Value = 3933.3333
$ws.Range("N32").Value = -4585.3333

$ws.Range("H40").Value = 4027.2856
$ws.Range("I40").Value = 3340
$ws.Range("J40").Value = 5745.5
$ws.Range("K40").Value = 3340
$ws.Range("L40").Value = 5745.5
$ws.Range("M40").Value = -3165
$ws.Range("N40").Value = -6095.5

$ws.Range("H51").Value = 25006480
$ws.Range("J51").Value = 9299.666999999999
$ws.Range("L51").Value = 9299.666999999999
$ws.Range("N51").Value = -10267.667

$ws.Range("H54").Value = 30833
$ws.Range("I54").Value = 20000
$ws.Range("J54").Value = 41666
$ws.Range("K54").Value = 20000
$ws.Range("L54").Value = 41666
$ws.Range("M54").Value = -19514
$ws.Range("N54").Value = -42638

$ws.Range("H62").Value = 4549.095
$ws.Range("I62").Value = 4281.676
$ws.Range("J62").Value = 6528
$ws.Range("K62").Value = 4281.676
$ws.Range("L62").Value = 6528
$ws.Range("M62").Value = -3657.676
$ws.Range("N62").Value = -7776

$ws.Range("H65").Value = 4549.095
$ws.Range("I65").Value = 4281.676
$ws.Range("J65").Value = 6528
$ws.Range("K65").Value = 21408.38
$ws.Range("L65").Value = 32640
$ws.Range("M65").Value = -18288.38
$ws.Range("N65").Value = -38880

$ws.Range("H88").Value = 7024.5835
$ws.Range("I88").Value = 2334.75
$ws.Range("J88").Value = 9369.5
$ws.Range("K88").Value = 2334.75
$ws.Range("L88").Value = 9369.5
$ws.Range("M88").Value = -1928.75
$ws.Range("N88").Value = -10181.5

$ws.Range("H91").Value = 7024.5835
$ws.Range("I91").Value = 2334.75
$ws.Range("J91").Value = 9369.5
$ws.Range("K91").Value = 2334.75
$ws.Range("L91").Value = 9369.5
$ws.Range("M91").Value = -930.75
$ws.Range("N91").Value = -12177.5

$ws.Range("H98").Value = 2414.0454
$ws.Range("I98").Value = 2578.7
$ws.Range("K98").Value = 2578.7
$ws.Range("M98").Value = -1080.7

$ws.Range("H122").Value = 2414.0454
$ws.Range("I122").Value = 2578.7
$ws.Range("K122").Value = 7736.099999999999
$ws.Range("M122").Value = -5286.099999999999

$ws.Range("H138").Value = 17333.734
$ws.Range("J138").Value = 94388.09
$ws.Range("L138").Value = 283164.27
$ws.Range("N138").Value = -293444.27

$ws.Range("H141").Value = 1507
$ws.Range("I141").Value = 1507
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4521
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = 659

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2371.1428
$ws.Range("I45").Value = 1679.85
$ws.Range("K45").Value = 1679.85
$ws.Range("M45").Value = -1302.85

$ws.Range("H61").Value = 10030.154
$ws.Range("I61").Value = 1311.75
$ws.Range("K61").Value = 1311.75
$ws.Range("M61").Value = -1099.75

$ws.Range("H122").Value = 2479.5557
$ws.Range("I122").Value = 2400.1428
$ws.Range("K122").Value = 7200.428400000001
$ws.Range("M122").Value = -4750.428400000001

$ws.Range("H132").Value = 1240.5186
$ws.Range("I132").Value = 939.76
$ws.Range("K132").Value = 2819.28
$ws.Range("M132").Value = -289.2799999999997

$ws.Range("H133").Value = 78000
$ws.Range("J133").Value = 78000
$ws.Range("L133").Value = 78000
$ws.Range("N133").Value = -83060

$ws.Range("H136").Value = 10030.154
$ws.Range("I136").Value = 1311.75
$ws.Range("K136").Value = 3935.25
$ws.Range("M136").Value = -1385.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 452.5
$ws.Range("I10").Value = 452.5
$ws.Range("K10").Value = 452.5
$ws.Range("M10").Value = -313.5

$ws.Range("H134").Value = 2598.2
$ws.Range("I134").Value = 1997.75
$ws.Range("K134").Value = 5993.25
$ws.Range("M134").Value = -3458.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3825.8
$ws.Range("I139").Value = 2282.25
$ws.Range("K139").Value = 6846.75
$ws.Range("M139").Value = -1706.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3856.7144
$ws.Range("I126").Value = 1999.75
$ws.Range("K126").Value = 5999.25
$ws.Range("M126").Value = -3529.25

$ws.Range("H132").Value = 4046.4138
$ws.Range("I132").Value = 3972.8333
$ws.Range("K132").Value = 11918.4999
$ws.Range("M132").Value = -9388.499899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 15299.889

$ws.Range("H25").Value = 50003.5
$ws.Range("I25").Value = 20007
$ws.Range("K25").Value = 20007
$ws.Range("M25").Value = -19777

$ws.Range("H46").Value = 4310.5454
$ws.Range("I46").Value = 800
$ws.Range("J46").Value = 5627
$ws.Range("K46").Value = 800
$ws.Range("L46").Value = 5627
$ws.Range("M46").Value = -612
$ws.Range("N46").Value = -6003

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 40062
$ws.Range("J49").Value = 40062
$ws.Range("L49").Value = 40062
$ws.Range("N49").Value = -40522

$ws.Range("H51").Value = 71000
$ws.Range("I51").Value = 51333.332
$ws.Range("K51").Value = 51333.332
$ws.Range("M51").Value = -50823.332

$ws.Range("H52").Value = 50047
$ws.Range("J52").Value = 50047
$ws.Range("L52").Value = 50047
$ws.Range("N52").Value = -50499

$ws.Range("H126").Value = 253184.5
$ws.Range("I126").Value = 2121.7693
$ws.Range("K126").Value = 6365.3079
$ws.Range("M126").Value = -3895.3079

$ws.Range("H132").Value = 2360.6
$ws.Range("I132").Value = 1479.7142
$ws.Range("K132").Value = 4439.142599999999
$ws.Range("M132").Value = -1909.142599999999

$ws.Range("H136").Value = 17384.648
$ws.Range("I136").Value = 18568.766
$ws.Range("K136").Value = 55706.298
$ws.Range("M136").Value = -53156.298

Write-Output "Applied all profit updates"
